$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (row 1): columns reordered / renamed ---
$ws.Range("B1").Value = "expression"
$ws.Range("C1").Value = "maximum"
$ws.Range("D1").Value = "minimum"
$ws.Range("E1").Value = "non_negative"
$ws.Range("F1").Value = "standard_error"
$ws.Range("G1").Value = "value"
$ws.Range("H1").Value = "vary"

# --- Row 2: pure_list.1 ---
$ws.Range("B2").Value = "None"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = $false
$ws.Range("F2").Value = "None"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = $true

# --- Row 3: pure_list.2 ---
$ws.Range("B3").Value = "None"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = $false
$ws.Range("F3").Value = "None"
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = $true

# --- Row 4: list_with_options.1 ---
$ws.Range("B4").Value = "None"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = "None"
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = $false

# --- Row 5: list_with_options.2 ---
$ws.Range("B5").Value = "None"
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = $false
$ws.Range("F5").Value = "None"
$ws.Range("G5").Value = 4
$ws.Range("H5").Value = $false

# --- Row 6: verbose_list.all_defaults ---
$ws.Range("B6").Value = "None"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = "None"
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = $true

# --- Row 7: verbose_list.no_defaults ---
$ws.Range("B7").Value = "None"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = -1
$ws.Range("E7").Value = $true
$ws.Range("F7").Value = "None"
$ws.Range("G7").Value = 6
$ws.Range("H7").Value = $false

# --- Row 8: verbose_list.expression_only ---
$ws.Range("B8").Value = '$verbose_list.all_defaults + $verbose_list.no_defaults'
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = "None"
$ws.Range("G8").Value = 11
$ws.Range("H8").Value = $false
